$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that only held "1285870 - Marcos Villela Barcza" (under the
# "Docentes responsáveis:" label) is removed; every row below it shifts up.
$ws.Rows("13").Delete()

# After the shift, re-point the remaining content cells so each label's
# value is now the value that used to belong to the row below it (the
# "Bibliografia:" row loses its value entirely).
$ws.Range("B10").Value = "1285870 - Marcos Villela Barcza"
$ws.Range("C10").Value = "1285870 - Marcos Villela Barcza"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Copy from the existing "01/01/2016" cell (Ativação:) so Excel keeps it
# as literal text instead of auto-converting the date-like string to a
# date serial number.
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

$ws.Range("B18").Value = "1285870 - Marcos Villela Barcza"
$ws.Range("C18").Value = "1285870 - Marcos Villela Barcza"

$ws.Range("B19").Value = "Reuniões periódicas com o orientador e realização do trabalho de conclusão de curso conforme orientação e apresentação de uma monografia final, conforme norma do Departamento de Engenharia Química."
$ws.Range("C19").Value = "Reuniões periódicas com o orientador e realização do trabalho de conclusão de curso conforme orientação e apresentação de uma monografia final, conforme norma do Departamento de Engenharia Química."

$ws.Range("B20").Value = "Avaliação da monografia perante uma banca examinadora composta por 3 (três) membros, obrigatoriamente docentes da Escola de Engenharia de Lorena (EEL)."
$ws.Range("C20").Value = "Avaliação da monografia perante uma banca examinadora composta por 3 (três) membros, obrigatoriamente docentes da Escola de Engenharia de Lorena (EEL)."

$ws.Range("B21").Value = "Reapresentação da monografia, preferencialmente para a mesma banca, com as modificações sugeridas para uma nova avaliação."
$ws.Range("C21").Value = "Reapresentação da monografia, preferencialmente para a mesma banca, com as modificações sugeridas para uma nova avaliação."
